$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two rows "RM 232" (orig row 26) and "SC 92" (orig row 28) were removed from
# the source data, which shifts every following row up. Delete the first one,
# then (at its now-shifted position) delete the second one.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# Remaining per-cell value edits (row numbers are POST-deletion / final positions).
$ws.Range("E5").Value = ""
$ws.Range("F6").Value = 16.43
$ws.Range("E8").Value = -6.6
$ws.Range("F11").Value = 17.65
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = ""
$ws.Range("E14").Value = -5.4
$ws.Range("F17").Value = ""
$ws.Range("E18").Value = ""
$ws.Range("F25").Value = 16.6
$ws.Range("B26").Value = -20.2
$ws.Range("B27").Value = ""
$ws.Range("F31").Value = ""
$ws.Range("F32").Value = ""
$ws.Range("D33").Value = -14.1
